$wb = $excel.ActiveWorkbook

# 1. Update the handoff status text everywhere it is used
#    ("Ready for handoff" -> "Handoff transform failed")
$newStatus = "Handoff transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = $newStatus

# 2. For both the zh-cn and de-de detail sheets, the previous handoff
#    attempt is rolled back: the "Latest Handoff File" link is removed,
#    the handoff datetime is reset and the handoff reason flips from
#    "Include" to "Ignored".
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $linksToRemove = @()
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$C$2') {
            $linksToRemove += $h
        }
    }
    foreach ($h in $linksToRemove) {
        $h.Delete()
    }

    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
